$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("2018 LEAVE CREDITS")
$ws2.Range("K93").Value = "10/9,11/2023"
